# NIT-9009155035.xlsx — "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Adds two new worker records (WARNER ENRIQUE AHUMADA GONZALEZ, periods 2011 and 2010)
# above the existing record (ANYI VANESSA RAMIREZ BEDOYA, period 2406) in the
# "Estado de Cuenta" table, and refreshes the summary totals accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the data table: turn the single data row (16) into three rows (16-18) ----
# Insert two fresh rows right after the existing data row, then fan the existing
# row's formatting down into them (Copy-to-destination keeps the same cell styles
# instead of minting brand-new ones for every paste).
$ws.Rows(17).Insert()
$ws.Rows(17).Insert()

$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# Row 18 keeps the pre-existing worker (ANYI VANESSA RAMIREZ BEDOYA / 2406) untouched.

# Row 16 -> new worker, period 2011
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73165726"
$ws.Range("D16").Value = "WARNER ENRIQUE AHUMADA GONZALEZ"
$ws.Range("E16").Value = "2011"
$ws.Range("F16").Value = 35112
$ws.Range("G16").Value = 877803

# Row 17 -> same worker, period 2010
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73165726"
$ws.Range("D17").Value = "WARNER ENRIQUE AHUMADA GONZALEZ"
$ws.Range("E17").Value = "2010"
$ws.Range("F17").Value = 35112
$ws.Range("G17").Value = 877803

# --- Refresh the summary block above the table -------------------------------------
$ws.Range("E11").Value = 122224   # VALOR MORA total (35112 + 35112 + 52000)
$ws.Range("C13").Value = 2        # Cant. Trabajadores
$ws.Range("F13").Value = 3        # Cant. Periodos
